# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before the
# "Late" column (column N), which pushes "Late" / "heading" (Original) /
# "Outstanding" one column to the right. Then make "Repayment schedule" the
# active sheet/tab with cell R9 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the "Late" column (N).
$ws.Columns("N").Insert() | Out-Null

# The newly inserted column inherits the width of the column to its left (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab with R9 selected.
$ws.Activate() | Out-Null
$ws.Range("R9").Select() | Out-Null
